$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Donte DiVincenzo",
    "Josh Hart",
    "Dyson Daniels",
    "Andrew Wiggins",
    "Chris Boucher",
    "Michael Porter Jr.",
    "De'Andre Hunter",
    "Alperen Sengün",
    "Victor Wembanyama",
    "Santi Aldama",
    "Kristaps Porzingis",
    "Donovan Mitchell",
    "Domantas Sabonis",
    "Deandre Ayton",
    "Malik Beasley",
    "Cam Thomas"
)

$positions = @(
    "PG,SG,SF",
    "SG,SF,PF",
    "PG,SG,SF",
    "SF,PF",
    "PF,C",
    "SF,PF",
    "SF,PF",
    "C",
    "C",
    "PF,C",
    "PF,C",
    "PG,SG",
    "C",
    "C",
    "SG,SF",
    "SG,SF"
)

$teams = @(
    "Minnesota Timberwolves",
    "New York Knicks",
    "Atlanta Hawks",
    "Golden State Warriors",
    "Toronto Raptors",
    "Denver Nuggets",
    "Atlanta Hawks",
    "Houston Rockets",
    "San Antonio Spurs",
    "Memphis Grizzlies",
    "Boston Celtics",
    "Cleveland Cavaliers",
    "Sacramento Kings",
    "Portland Trail Blazers",
    "Detroit Pistons",
    "Brooklyn Nets"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
    $ws.Cells.Item($i + 2, 2).Value = $positions[$i]
    $ws.Cells.Item($i + 2, 3).Value = $teams[$i]
}

# The old table had one extra (17th) data row; remove it so the sheet ends
# at row 17 (header + 16 players), matching the updated roster.
$ws.Rows(18).Delete()
